$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in rows 2-5 to the new cluster order/counts
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 342

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 178

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 157

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 37

# Row 6 is no longer part of the data - delete it entirely
$ws.Rows.Item(6).Delete()
